$d = $word.ActiveDocument

# 1. Remove the paragraph that only contains the horizontal-rule picture
#    (the <w:p> right after the "DOC-to-Markdown" heading, before "Disclaimer").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $cand.Range.Delete()
        break
    }
}

# 2. Re-split the disclaimer paragraph's two long runs into several runs that
#    break on (roughly) 79-character boundaries, with standalone single-space
#    runs between them, matching how the markdown-to-docx converter wrapped
#    the source text. The quoted '‘as is’' runs in the middle of the
#    paragraph are left untouched.
$disclaimer = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("This repository is a scientific product")) {
        $disclaimer = $cand
        break
    }
}
$bodyRange = $d.Range($disclaimer.Range.Start, $disclaimer.Range.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This repository is a scientific product and is not official</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">communication of the National Oceanic and Atmospheric Administration, or</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">the United States Department of Commerce. All NOAA GitHub project</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">content is provided on an</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">‘</w:t></w:r><w:r><w:t xml:space="preserve">as is</w:t></w:r><w:r><w:t xml:space="preserve">’</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">basis and the user assumes</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">responsibility for its use. Any claims against the Department of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Commerce or Department of Commerce bureaus stemming from the use of this</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">GitHub project will be governed by all applicable Federal law. Any</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">reference to specific commercial products, processes, or services by</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">service mark, trademark, manufacturer, or otherwise, does not constitute</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">or imply their endorsement, recommendation or favoring by the Department</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">of Commerce. The Department of Commerce seal and logo, or the seal and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">logo of a DOC bureau, shall not be used in any manner to imply</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">endorsement of any commercial product or activity by DOC or the United</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">States Government.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$bodyRange.InsertXML($xml)
